$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 76
$ws.Range("H76").Value = 3637.8572
$ws.Range("I76").Value = 3555
$ws.Range("J76").Value = 3700
$ws.Range("K76").Value = 3555
$ws.Range("L76").Value = 3700
$ws.Range("M76").Value = -3240
$ws.Range("N76").Value = -4330
# row 79
$ws.Range("H79").Value = 3637.8572
$ws.Range("I79").Value = 3555
$ws.Range("J79").Value = 3700
$ws.Range("K79").Value = 3555
$ws.Range("L79").Value = 3700
$ws.Range("M79").Value = -2463
$ws.Range("N79").Value = -5884
# row 100
$ws.Range("H100").Value = 66668052
$ws.Range("I100").Value = 100000630
$ws.Range("J100").Value = 2901
$ws.Range("K100").Value = 100000630
$ws.Range("L100").Value = 2901
$ws.Range("M100").Value = -100000089
$ws.Range("N100").Value = -3983
# row 113
$ws.Range("H113").Value = 34486630
$ws.Range("I113").Value = 83336670
$ws.Range("J113").Value = 4249.7646
$ws.Range("K113").Value = 83336670
$ws.Range("L113").Value = 4249.7646
$ws.Range("M113").Value = -83333416
$ws.Range("N113").Value = -10757.7646

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 102
$ws.Range("H102").Value = 1097.4
$ws.Range("I102").Value = 871.75
$ws.Range("K102").Value = 871.75
$ws.Range("M102").Value = 750.25
# row 110
$ws.Range("H110").Value = 526
$ws.Range("I110").Value = 491.75
$ws.Range("K110").Value = 491.75
$ws.Range("M110").Value = 1553.25

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 1634.8438
$ws.Range("I86").Value = 1403.8572
$ws.Range("J86").Value = 3251.75
$ws.Range("K86").Value = 1403.8572
$ws.Range("L86").Value = 3251.75
$ws.Range("M86").Value = -280.8571999999999
$ws.Range("N86").Value = -5497.75
# row 89
$ws.Range("H89").Value = 1634.8438
$ws.Range("I89").Value = 1403.8572
$ws.Range("J89").Value = 3251.75
$ws.Range("K89").Value = 7019.286
$ws.Range("L89").Value = 16258.75
$ws.Range("M89").Value = -1403.286
$ws.Range("N89").Value = -27490.75
# row 99
$ws.Range("H99").Value = 1000.1212
$ws.Range("I99").Value = 796.0741
$ws.Range("J99").Value = 1918.3334
$ws.Range("K99").Value = 796.0741
$ws.Range("L99").Value = 1918.3334
$ws.Range("M99").Value = 701.9259
$ws.Range("N99").Value = -4914.3334
# row 107
$ws.Range("H107").Value = 5654.3335
$ws.Range("I107").Value = 4975
$ws.Range("K107").Value = 4975
$ws.Range("M107").Value = -3055
# row 134
$ws.Range("H134").Value = 3963.6333
$ws.Range("I134").Value = 3963.6333
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11890.8999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -9355.8999
$ws.Range("N134").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 1462
$ws.Range("I16").Value = 1596
$ws.Range("J16").Value = 1294.5
$ws.Range("K16").Value = 1596
$ws.Range("L16").Value = 1294.5
$ws.Range("M16").Value = -1309
$ws.Range("N16").Value = -1868.5
# row 69
$ws.Range("H69").Value = 8891.200000000001
$ws.Range("I69").Value = 8441
$ws.Range("K69").Value = 8441
$ws.Range("M69").Value = -7692
# row 72
$ws.Range("H72").Value = 8891.200000000001
$ws.Range("I72").Value = 8441
$ws.Range("K72").Value = 25323
$ws.Range("M72").Value = -21579
# row 107
$ws.Range("H107").Value = 1246.8422
$ws.Range("I107").Value = 754
$ws.Range("J107").Value = 1794.4445
$ws.Range("K107").Value = 754
$ws.Range("L107").Value = 1794.4445
$ws.Range("M107").Value = 1166
$ws.Range("N107").Value = -5634.4445
# row 113
$ws.Range("H113").Value = 1462
$ws.Range("I113").Value = 1596
$ws.Range("J113").Value = 1294.5
$ws.Range("K113").Value = 1596
$ws.Range("L113").Value = 1294.5
$ws.Range("M113").Value = 574
$ws.Range("N113").Value = -5634.5
# row 132
$ws.Range("H132").Value = 17655.176
$ws.Range("I132").Value = 24015.088
$ws.Range("K132").Value = 72045.264
$ws.Range("M132").Value = -69515.264
# row 134
$ws.Range("H134").Value = 1032.4
$ws.Range("I134").Value = 1032.4
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3097.2
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -562.2000000000003
$ws.Range("N134").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 131
$ws.Range("H131").Value = 789.11
$ws.Range("J131").Value = 794.0505000000001
$ws.Range("L131").Value = 2382.1515
$ws.Range("N131").Value = -12462.1515

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 62
$ws.Range("H62").Value = 20085
$ws.Range("J62").Value = 20085
$ws.Range("L62").Value = 20085
$ws.Range("N62").Value = -21457
# row 65
$ws.Range("H65").Value = 20085
$ws.Range("J65").Value = 20085
$ws.Range("L65").Value = 60255
$ws.Range("N65").Value = -67119
# row 107
$ws.Range("H107").Value = 5917790.5
$ws.Range("I107").Value = 244.44444
$ws.Range("J107").Value = 19232268
$ws.Range("K107").Value = 244.44444
$ws.Range("L107").Value = 19232268
$ws.Range("M107").Value = 1675.55556
$ws.Range("N107").Value = -19236108
# row 123
$ws.Range("H123").Value = 6977.727
$ws.Range("I123").Value = 2960
$ws.Range("J123").Value = 10325.833
$ws.Range("K123").Value = 2960
$ws.Range("L123").Value = 10325.833
$ws.Range("M123").Value = -510
$ws.Range("N123").Value = -15225.833

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 68
$ws.Range("H68").Value = 2106.923
$ws.Range("I68").Value = 1860
$ws.Range("J68").Value = 2261.25
$ws.Range("K68").Value = 1860
$ws.Range("L68").Value = 2261.25
$ws.Range("M68").Value = -1111
$ws.Range("N68").Value = -3759.25
# row 71
$ws.Range("H71").Value = 2106.923
$ws.Range("I71").Value = 1860
$ws.Range("J71").Value = 2261.25
$ws.Range("K71").Value = 9300
$ws.Range("L71").Value = 11306.25
$ws.Range("M71").Value = -5556
$ws.Range("N71").Value = -18794.25
# row 132
$ws.Range("H132").Value = 2290.3809
$ws.Range("I132").Value = 1600.1818
$ws.Range("J132").Value = 3049.6
$ws.Range("K132").Value = 4800.5454
$ws.Range("L132").Value = 9148.799999999999
$ws.Range("M132").Value = -2270.5454
$ws.Range("N132").Value = -14208.8

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 81
$ws.Range("H81").Value = 125001730
$ws.Range("I81").Value = 1971.5714
$ws.Range("J81").Value = 1000000000
$ws.Range("K81").Value = 3943.1428
$ws.Range("L81").Value = 2000000000
$ws.Range("M81").Value = -2882.1428
$ws.Range("N81").Value = -2000002122
# row 84
$ws.Range("H84").Value = 125001730
$ws.Range("I84").Value = 1971.5714
$ws.Range("J84").Value = 1000000000
$ws.Range("K84").Value = 19715.714
$ws.Range("L84").Value = 10000000000
$ws.Range("M84").Value = -14411.714
$ws.Range("N84").Value = -10000010608
# row 100
$ws.Range("H100").Value = 442.2857
$ws.Range("I100").Value = 419.4
$ws.Range("K100").Value = 838.8
$ws.Range("M100").Value = -297.8
# row 132
$ws.Range("H132").Value = 1018.119
$ws.Range("I132").Value = 766.0357
$ws.Range("J132").Value = 1522.2858
$ws.Range("K132").Value = 2298.1071
$ws.Range("L132").Value = 4566.857400000001
$ws.Range("M132").Value = 231.8928999999998
$ws.Range("N132").Value = -9626.857400000001

Write-Host "Edit complete"